# Reorder the player roster rows (rows 2-19) on Sheet1.
# Each player keeps their original Position/Team pairing; only the
# row order in which the players are listed changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("Austin Reaves",       "PG,SG", "Los Angeles Lakers"),
    @("Tyrese Haliburton",   "PG,SG", "Indiana Pacers"),
    @("Stephen Curry",       "PG,SG", "Golden State Warriors"),
    @("Darius Garland",      "PG",    "Cleveland Cavaliers"),
    @("Keegan Murray",       "SF,PF", "Sacramento Kings"),
    @("Daniel Gafford",      "PF,C",  "Dallas Mavericks"),
    @("Karl-Anthony Towns",  "PF,C",  "New York Knicks"),
    @("Kevin Durant",        "SF,PF", "Phoenix Suns"),
    @("Jarrett Allen",       "C",     "Cleveland Cavaliers"),
    @("Jalen Duren",         "C",     "Detroit Pistons"),
    @("Trey Murphy III",     "SF,PF", "New Orleans Pelicans"),
    @("Tyrese Maxey",        "PG,SG", "Philadelphia 76ers"),
    @("Jrue Holiday",        "PG,SG", "Boston Celtics"),
    @("OG Anunoby",          "SF,PF", "New York Knicks"),
    @("Mark Williams",       "C",     "Charlotte Hornets"),
    @("Franz Wagner",        "SF,PF", "Orlando Magic"),
    @("Jalen Johnson",       "SF,PF", "Atlanta Hawks"),
    @("Dereck Lively II",    "C",     "Dallas Mavericks")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
